$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster

Write-Host "BEFORE sm colors:"
$cs = $sm.ColorScheme
for ($i=1; $i -le $cs.Count; $i++) {
  Write-Host " $i : $($cs.Colors($i).RGB)"
}

$sm.ColorScheme = $nm.ColorScheme

Write-Host "AFTER sm colors:"
$cs2 = $sm.ColorScheme
for ($i=1; $i -le $cs2.Count; $i++) {
  Write-Host " $i : $($cs2.Colors($i).RGB)"
}
